# Update gh-pages to output generated at 456a3b4
# Applies refreshed "want to go" counts (col F) and two swapped cover images (col I)
# across the 展览/演出/本地生活/全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1283
$ws.Range("F5").Value = 14
$ws.Range("F6").Value = 6852
$ws.Range("F7").Value = 1816
$ws.Range("F8").Value = 6402
$ws.Range("F9").Value = 145
$ws.Range("F10").Value = 1956
$ws.Range("F11").Value = 523
$ws.Range("F12").Value = 20
$ws.Range("F14").Value = 38
$ws.Range("F18").Value = 8058
$ws.Range("F19").Value = 142
$ws.Range("F23").Value = 1757
$ws.Range("F30").Value = 2
$ws.Range("F31").Value = 1837
$ws.Range("F32").Value = 818
$ws.Range("F33").Value = 394
$ws.Range("F36").Value = 19
$ws.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202409/LU32zDTR1725617506119.jpeg"
$ws.Range("F38").Value = 89
$ws.Range("F39").Value = 3930
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 374
$ws.Range("F14").Value = 11
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9548
$ws.Range("F4").Value = 694
$ws.Range("F5").Value = 282
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9548
$ws.Range("F4").Value = 694
$ws.Range("F5").Value = 1283
$ws.Range("F9").Value = 374
$ws.Range("F10").Value = 6852
$ws.Range("F11").Value = 282
$ws.Range("F12").Value = 1816
$ws.Range("F13").Value = 6402
$ws.Range("F14").Value = 145
$ws.Range("F15").Value = 1956
$ws.Range("F18").Value = 523
$ws.Range("F19").Value = 38
$ws.Range("F21").Value = 11
$ws.Range("F24").Value = 8058
$ws.Range("F25").Value = 142
$ws.Range("F28").Value = 1757
$ws.Range("F32").Value = 1837
$ws.Range("F33").Value = 818
$ws.Range("F39").Value = 19
$ws.Range("I39").Value = "//i2.hdslb.com/bfs/openplatform/202409/LU32zDTR1725617506119.jpeg"
$ws.Range("F42").Value = 89
$ws.Range("F44").Value = 3930
